# quantitative evaluation plots and execution metrics analysis
# -------------------------------------------------------------
# The underlying per-test-case categorical results for configuration "11"
# (compilation_success / execution_without_error / test_pass columns B/C/D)
# were corrected; this ripples into the derived summary tables
# (Categorical_Summary, Summary_All_Configs and Success_Rates row 12, which
# all summarise configuration "11") that were stored as plain cached values
# rather than live formulas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Sheet "11" - raw per-test-case Yes/No/Unknown categorical data
# ---------------------------------------------------------------
$ws11 = $wb.Worksheets("11")

$ws11.Range("D4").Value = "No"

$ws11.Range("B6").Value = "Yes"
$ws11.Range("C6").Value = "No"

$ws11.Range("B8").Value = "Yes"
$ws11.Range("C8").Value = "No"

$ws11.Range("C9").Value = "No"
$ws11.Range("D9").Value = "No"

$ws11.Range("B10").Value = "No"
$ws11.Range("C10").Value = "Unknown"
$ws11.Range("D10").Value = "Unknown"

$ws11.Range("B11").Value = "Yes"
$ws11.Range("C11").Value = "No"

$ws11.Range("B12").Value = "Yes"
$ws11.Range("C12").Value = "No"

$ws11.Range("B13").Value = "Yes"
$ws11.Range("C13").Value = "No"

$ws11.Range("B14").Value = "Yes"
$ws11.Range("C14").Value = "No"

$ws11.Range("B15").Value = "Yes"
$ws11.Range("C15").Value = "No"

$ws11.Range("B16").Value = "Yes"
$ws11.Range("C16").Value = "No"

$ws11.Range("B19").Value = "Yes"
$ws11.Range("C19").Value = "Yes"
$ws11.Range("D19").Value = "No"

# ---------------------------------------------------------------
# 2) Categorical_Summary - row 12 (Configuration "11") recomputed counts
# ---------------------------------------------------------------
$wsCat = $wb.Worksheets("Categorical_Summary")

$wsCat.Range("B12").Value = 16
$wsCat.Range("C12").Value = 88.88888888888889
$wsCat.Range("D12").Value = 2
$wsCat.Range("E12").Value = 11.11111111111111
$wsCat.Range("F12").Value = 88.88888888888889
$wsCat.Range("G12").Value = 16
$wsCat.Range("H12").Value = 2

$wsCat.Range("J12").Value = 7
$wsCat.Range("K12").Value = 38.88888888888889
$wsCat.Range("L12").Value = 2
$wsCat.Range("M12").Value = 11.11111111111111
$wsCat.Range("N12").Value = 9
$wsCat.Range("O12").Value = 50
$wsCat.Range("P12").Value = 38.88888888888889
$wsCat.Range("Q12").Value = 7
$wsCat.Range("R12").Value = 9

$wsCat.Range("V12").Value = 3
$wsCat.Range("W12").Value = 16.66666666666666
$wsCat.Range("X12").Value = 5
$wsCat.Range("Y12").Value = 27.77777777777778
$wsCat.Range("Z12").Value = 27.77777777777778
$wsCat.Range("AA12").Value = 5
$wsCat.Range("AB12").Value = 3

# ---------------------------------------------------------------
# 3) Success_Rates - row 12 (Configuration "11")
# ---------------------------------------------------------------
$wsRates = $wb.Worksheets("Success_Rates")

$wsRates.Range("B12").Value = 88.88888888888889
$wsRates.Range("C12").Value = 38.88888888888889
$wsRates.Range("D12").Value = 27.77777777777778

# ---------------------------------------------------------------
# 4) Summary_All_Configs - row 12 (Configuration "11"), columns N..AN
# ---------------------------------------------------------------
$wsSum = $wb.Worksheets("Summary_All_Configs")

$wsSum.Range("N12").Value = 16
$wsSum.Range("O12").Value = 88.88888888888889
$wsSum.Range("P12").Value = 2
$wsSum.Range("Q12").Value = 11.11111111111111
$wsSum.Range("R12").Value = 88.88888888888889
$wsSum.Range("S12").Value = 16
$wsSum.Range("T12").Value = 2

$wsSum.Range("V12").Value = 7
$wsSum.Range("W12").Value = 38.88888888888889
$wsSum.Range("X12").Value = 2
$wsSum.Range("Y12").Value = 11.11111111111111
$wsSum.Range("Z12").Value = 9
$wsSum.Range("AA12").Value = 50
$wsSum.Range("AB12").Value = 38.88888888888889
$wsSum.Range("AC12").Value = 7
$wsSum.Range("AD12").Value = 9

$wsSum.Range("AH12").Value = 3
$wsSum.Range("AI12").Value = 16.66666666666666
$wsSum.Range("AJ12").Value = 5
$wsSum.Range("AK12").Value = 27.77777777777778
$wsSum.Range("AL12").Value = 27.77777777777778
$wsSum.Range("AM12").Value = 5
$wsSum.Range("AN12").Value = 3
